# Update cryptocurrency price/volume figures per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D2:E51 are stored as text (inline strings), not numbers, so the
# NumberFormat is forced to text before the write and cleared back to the
# sheet's default ('General') afterwards -- this keeps values like
# '234.76' or '37.842.99' from being auto-converted into numeric cells.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.842.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E10").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.396.61"
$ws.Range("D12").ClearFormats()

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.75"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.43"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.772"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.091.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.766.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.69"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.87"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.54%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.71"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.52"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.91%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0991"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.60"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.463.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.64%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.05%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.38%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.68"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.33%  "
$ws.Range("E51").ClearFormats()

